$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 101, shifting existing rows 101-102 down to 102-103
$ws.Rows.Item(101).Insert()

# Fill in the new row 101 with the inserted record's data
$ws.Cells.Item(101, 1).Value = 11
$ws.Cells.Item(101, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(101, 3).Value = "Bíobío"
$ws.Cells.Item(101, 4).Value = 44628
$ws.Cells.Item(101, 5).Value = 8
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100108
$ws.Cells.Item(101, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(101, 9).Value = 100108002
$ws.Cells.Item(101, 10).Value = "Mango"
$ws.Cells.Item(101, 11).Value = "Sin especificar"
$ws.Cells.Item(101, 12).Value = "Primera"
$ws.Cells.Item(101, 13).Value = 200
$ws.Cells.Item(101, 14).Value = 7500
$ws.Cells.Item(101, 15).Value = 8000
$ws.Cells.Item(101, 16).Value = 7750
$ws.Cells.Item(101, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(101, 18).Value = "Ecuador"
$ws.Cells.Item(101, 19).Value = 1938
$ws.Cells.Item(101, 20).Value = 4
